$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the data_status column (B) for existing teams rows 31-34
$ws.Range("B31").Value = "pulled"
$ws.Range("B32").Value = "pulled"
$ws.Range("B33").Value = "pulled"
$ws.Range("B34").Value = "pulled"

# Add new row 40 first (West Virginia / pulled) so that "West Virginia"
# is interned into the shared-string table before "not-pulled" is used.
$ws.Range("A40").Value = "West Virginia"
$ws.Range("B40").Value = "pulled"

# Rows 35-36 marked not-pulled
$ws.Range("B35").Value = "not-pulled"
$ws.Range("B36").Value = "not-pulled"

$ws.Range("B37").Value = "pulled"

# Swap order: UCSB moves to row 38, Louisiana Tech moves to row 39
$ws.Range("A38").Value = "UCSB"
$ws.Range("B38").Value = "not-pulled"
$ws.Range("A39").Value = "Louisiana Tech"
$ws.Range("B39").Value = "not-pulled"

# Update view: clear topLeftCell scroll, select A41
$ws.Range("A41").Select()
